# Generate Report for Handback
# Adds a new handback record (b013591f-7584-4a73-a059-a39a2aac7c99) as a
# new row (row 3) on each of the three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$newFile      = "b013591f-7584-4a73-a059-a39a2aac7c99.md"
$newFilePath  = "e2e\b013591f-7584-4a73-a059-a39a2aac7c99.md"
$genDate      = "2017-02-09 08:01:57"

$zhXlf        = "b013591f-7584-4a73-a059-a39a2aac7c99.408e5bacc2e22a8dd0b7d42bca0f848975663a38.zh-cn.xlf"
$zhHoDate     = "2017-02-09 08:01:40"
$zhHbDate     = "2017-02-09 08:02:39"

$deXlf        = "b013591f-7584-4a73-a059-a39a2aac7c99.408e5bacc2e22a8dd0b7d42bca0f848975663a38.de-de.xlf"
$deHoDate     = "2017-02-09 08:01:57"
$deHbDate     = "2017-02-09 08:03:03"

# ---------------------------------------------------------------------
# Sheet "Overview" - append row 3 to the Overview table
# ---------------------------------------------------------------------
$tOverview = $ws1.ListObjects.Item(1)
$tOverview.ListRows.Add() | Out-Null

$ws1.Range("A3").Value = $newFile
$ws1.Range("C3").Value = ".md"
$ws1.Range("E3").Value = "Handed back: in sync with en-US"
$ws1.Range("F3").Value = "Handed back: in sync with en-US"
$ws1.Range("G3").Value = $genDate
$ws1.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a28272383e24ad73c8d92fe07e2d8956f1af8d3/e2e/b013591f-7584-4a73-a059-a39a2aac7c99.md", "", "", $newFilePath) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn" - append row 3 to the zh-cn table
# ---------------------------------------------------------------------
$tZh = $ws2.ListObjects.Item(1)
$tZh.ListRows.Add() | Out-Null

$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "'True"
$ws2.Range("G3").Value = $zhXlf
$ws2.Range("H3").Value = $zhHoDate
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("K3").Value = $zhXlf
$ws2.Range("L3").Value = $zhHbDate
$ws2.Range("L3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("O3").Value = "'True"
$ws2.Range("Q3").Value = "'False"

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4548cc233b7534afa9f4d73e7b6839a56d1ddc01/e2e/b013591f-7584-4a73-a059-a39a2aac7c99.md", "", "", $newFile) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4548cc233b7534afa9f4d73e7b6839a56d1ddc01/e2e/b013591f-7584-4a73-a059-a39a2aac7c99.md", "", "", $newFile) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de" - append row 3 to the de-de table
# ---------------------------------------------------------------------
$tDe = $ws3.ListObjects.Item(1)
$tDe.ListRows.Add() | Out-Null

$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "'True"
$ws3.Range("G3").Value = $deXlf
$ws3.Range("H3").Value = $deHoDate
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("K3").Value = $deXlf
$ws3.Range("L3").Value = $deHbDate
$ws3.Range("L3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("O3").Value = "'True"
$ws3.Range("Q3").Value = "'False"

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/629838ff0f85d56630f06315559779a76c260f11/e2e/b013591f-7584-4a73-a059-a39a2aac7c99.md", "", "", $newFile) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/629838ff0f85d56630f06315559779a76c260f11/e2e/b013591f-7584-4a73-a059-a39a2aac7c99.md", "", "", $newFile) | Out-Null

Write-Output "Handback report row added to Overview, zh-cn and de-de sheets."
